# Update handback status timestamps (Generate Report for Handback)
$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-02-24 09:42:18"
$wsZhCn.Range("G4").Value = "2016-02-24 09:43:05"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-02-24 09:42:31"
$wsDeDe.Range("G4").Value = "2016-02-24 09:43:33"
